$d = $word.ActiveDocument

# --- 1 & 2: remove the " ***" suffix runs after the two "Dmg Die" enhancement
#            lines in the Demon Claws row. A single Replace-All Find/Execute
#            over the whole document content removes both occurrences.
$d.Content.Find.Execute(" ***", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- 3: fill in the previously-empty "Enhancements" cell for the Scales
#        row (table row 10, column 10) with "Armor / x3 / -- / 10P",
#        matching the formatting used throughout the table (Abadi, 7pt).
$table = $d.Tables.Item(1)
$cell = $table.Cell(10, 10)

$r1 = $cell.Range
$r1.Text = "Armor"
$r1.Font.Name = "Abadi"
$r1.Font.Size = 7
$r1.Font.SizeBi = 7

$cell = $table.Cell(10, 10)
$r2 = $cell.Range
$r2.Collapse(0)
$r2.MoveEnd(1, -1)
$r2.InsertAfter(" / x3 / -- / ")
$r2.Font.Name = "Abadi"
$r2.Font.Size = 7
$r2.Font.SizeBi = 7

$cell = $table.Cell(10, 10)
$r3 = $cell.Range
$r3.Collapse(0)
$r3.MoveEnd(1, -1)
$r3.InsertAfter("10P")
$r3.Font.Name = "Abadi"
$r3.Font.Size = 7
$r3.Font.SizeBi = 7

# Leave a "_GoBack" bookmark at the end of the cell's content, as the
# authoring session would.
$cell = $table.Cell(10, 10)
$endRange = $cell.Range
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRange)
